$d = $word.ActiveDocument

# Around the "_GoBack" bookmark we have:
#   ...<w:r>Step 4 : </w:r><w:bookmarkStart .../><w:bookmarkEnd .../><w:r>Pushing project to remot repository :</w:r>
# The author retyped "Pushing project to remote" in front of the bookmark
# (fixing the "remot" -> "remote" typo) and trimmed the same text back out
# of the run that follows the bookmark, leaving " repository :" there.

$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Range.Start

$prefix = "Pushing project to remot"
$fixedChar = "e"

# Insert the corrected text immediately before the bookmark.
$insertRange = $d.Range($bmStart, $bmStart)
$insertRange.Text = $prefix
$insertRange.Font.Bold = 1
$insertRange.Font.Size = 14

$bm2 = $d.Bookmarks.Item("_GoBack")
$bmStart2 = $bm2.Range.Start

$insertRange2 = $d.Range($bmStart2, $bmStart2)
$insertRange2.Text = $fixedChar
$insertRange2.Font.Bold = 1
$insertRange2.Font.Size = 14

# Now remove the duplicated "Pushing project to remot" text from the start
# of the run that follows the bookmark, leaving " repository :" behind.
$bm3 = $d.Bookmarks.Item("_GoBack")
$bmEnd3 = $bm3.Range.End
$oldPrefixRange = $d.Range($bmEnd3, $bmEnd3 + $prefix.Length)
if ($oldPrefixRange.Text -eq $prefix) {
    $oldPrefixRange.Text = ""
}
